$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new blank rows above the current row 16 ("implicit masking" row),
# shifting it (and everything below) down by two rows.
$ws.Rows("16:17").Insert()

# --- New row 16: options_mrf_weighting -------------------------------------
$ws.Range("A16").Value = "options_mrf_weighting"
$ws.Range("B16").Value = "markov random field weighting"
$ws.Range("C16").Value = "When tissue class images are written out, a few iterations of a simple Markov Random Field`n(MRF) cleanup procedure are run. This parameter controls the strength of the MRF. Setting the`nvalue to zero will disable the cleanup."
$ws.Range("D16").Value = "number"
$ws.Range("E16").Value = "float"
$ws.Range("F16").Value = 1
$ws.Range("I16").Value = "segmentation"

# --- New row 17: options_cleanup -------------------------------------------
$ws.Range("A17").Value = "options_cleanup"
$ws.Range("B17").Value = "clean up"
$ws.Range("C17").Value = "This uses a crude routine for extracting the brain from segmented images.`nIt begins by taking the white matter, and eroding it a couple of times to get rid of any odd`nvoxels. The algorithm continues on to do conditional dilations for several iterations, where the`ncondition is based upon gray or white matter being present.This identified region is then used to`nclean up the grey and white matter partitions. Note that the fluid class will also be cleaned, such`nthat aqueous and vitreous humour in the eyeballs, as well as other assorted fluid regions (except`nCSF) will be removed.`nIf you find pieces of brain being chopped out in your data, then you may wish to disable or`ntone down the cleanup procedure. Note that the procedure uses a number of assumptions about`nwhat each tissue class refers to. If a different set of tissue priors are used, then this routine should`nbe disabled.0=None,1=light,2=thorough"
$ws.Range("D17").Value = "number"
$ws.Range("E17").Value = "int"
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = "[0,1,2]"
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = "segmentation"

# Make sure every populated cell in the two new rows carries the same
# wrap-text format ("style 2") used throughout the rest of the sheet.
$ws.Range("A16:F16").WrapText = $true
$ws.Range("I16").WrapText = $true
$ws.Range("A17:I17").WrapText = $true

# Row heights grew to fit the new (longer) wrapped text.
$ws.Rows("16").RowHeight = 136
$ws.Rows("17").RowHeight = 409.6

# Put the selection on the newly added "cleanup" row, matching the saved view.
$ws.Range("A17").Select()
